$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: new "t1_omics_fields" entry ---
$ws.Cells.Item(2, 1).Value = "t1_omics_fields"
$ws.Cells.Item(2, 2).Value = "Omic fields targeted by xOmics"
$ws.Cells.Item(2, 2).WrapText = $true
$ws.Cells.Item(2, 3).Clear()

# --- Row 3: new "t2_quantification_methods" entry ---
$ws.Cells.Item(3, 1).Value = "t2_quantification_methods"
$ws.Cells.Item(3, 2).Value = "Quantification methods used in omic fields"
$ws.Cells.Item(3, 3).Clear()

# --- Row 4: renamed from t1_overview_datasets -> t3_overview_datasets ---
$ws.Cells.Item(4, 1).Value = "t3_overview_datasets"
$ws.Cells.Item(4, 2).Value = "Omics example datasets"
$ws.Cells.Item(4, 3).Value = "aa.load_dataset"

# --- Row 5: renamed from t2_omics_analysis_tools -> t4_omics_analysis_tools ---
$ws.Cells.Item(5, 1).Value = "t4_omics_analysis_tools"
$ws.Cells.Item(5, 2).Value = "Analysis tools for omics data"

# --- Row 6: renamed from t3_omics_post_analysis_tools -> t5_omics_post_analysis_tools ---
$ws.Cells.Item(6, 1).Value = "t5_omics_post_analysis_tools"
$ws.Cells.Item(6, 2).Value = "Post-analysis tools for omics data"

# --- Row 7: renamed from t4_gene_enrichment_tools -> t6_enrichment_tools, text changed ---
$ws.Cells.Item(7, 1).Value = "t6_enrichment_tools"
$ws.Cells.Item(7, 2).Value = "Enrichment analysis tools"

# --- View: update selection to B9:B10 ---
$ws.Range("B9:B10").Select()

# --- Sheet default column width ---
$ws.StandardWidth = 8.875
